# 2019-12-05 Report.xlsx — add a small "Project / Start Date / End Date"
# header block (and an "Issues" / "level" / "count" sub-header) above the
# existing Issues-By-Priority-And-Count table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Issues By Priority And Count")

# Push the existing 9 data rows (Critical..Low) down by 5 rows so they land
# on rows 6-14, opening up rows 1-5 for the new header block (row 3 is left
# blank on purpose, matching the target layout).
$ws.Range("A1:A5").EntireRow.Insert()

# Row 1: column headers
$ws.Range("A1").Value = "Project"
$ws.Range("B1").Value = "Start Date"
$ws.Range("C1").Value = "End Date"

# Row 2: project name + date range
$ws.Range("A2").Value = "DEVTST"
$ws.Range("B2").Value = 43774
$ws.Range("C2").Value = 43804

# Format the two date cells as short dates (built-in numFmtId 14) and make
# sure both cells share the exact same style entry by formatting one cell
# then copying its format onto the other (avoids creating two duplicate
# style entries).
$ws.Range("B2").NumberFormat = "mm-dd-yy"
$ws.Range("B2").Copy()
$ws.Range("C2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 3 intentionally left blank.

# Row 4: section title
$ws.Range("A4").Value = "Issues"

# Row 5: table header for the priority/count breakdown below
$ws.Range("A5").Value = "level"
$ws.Range("B5").Value = "count"
